$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E columns to Text format so numeric/percent-looking strings
# are preserved verbatim as text (matching original inlineStr cells),
# instead of Excel auto-converting them into numbers/percentages.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = "328.23"
$ws.Range("E2").Value = "6.60%"
$ws.Range("D3").Value = "39.92"
$ws.Range("E3").Value = "8.37%"
$ws.Range("E4").Value = "11.66%"
$ws.Range("D5").Value = "0.08127"
$ws.Range("E5").Value = "4.73%"
$ws.Range("D6").Value = "4.563"
$ws.Range("E6").Value = "3.56%"
$ws.Range("D7").Value = "8.686"
$ws.Range("E7").Value = "4.40%"
$ws.Range("D8").Value = "1.962"
$ws.Range("E8").Value = "5.88%"
$ws.Range("E9").Value = "1.67%"
$ws.Range("D10").Value = "0.9424"
$ws.Range("E10").Value = "1.92%"
$ws.Range("D11").Value = "0.1292"
$ws.Range("E11").Value = "15.63%"
$ws.Range("D12").Value = "0.1993"
$ws.Range("E12").Value = "6.46%"
$ws.Range("D13").Value = "0.09193"
$ws.Range("E13").Value = "4.43%"
$ws.Range("D14").Value = "0.03474"
$ws.Range("E14").Value = "5.68%"
$ws.Range("D15").Value = "0.09615"
$ws.Range("E15").Value = "0.65%"
$ws.Range("D16").Value = "0.001308"
$ws.Range("E16").Value = "-6.00%"
$ws.Range("D17").Value = "0.006103"
$ws.Range("E17").Value = "-1.77%"
$ws.Range("D18").Value = "3.368"
$ws.Range("E18").Value = "-0.73%"
$ws.Range("D19").Value = "0.3535"
$ws.Range("E19").Value = "2.37%"
$ws.Range("D20").Value = "7.568"
$ws.Range("E20").Value = "18.93%"
$ws.Range("E21").Value = "8.95%"
$ws.Range("D22").Value = "0.2421"
$ws.Range("E22").Value = "3.53%"
$ws.Range("D23").Value = "0.04443"
$ws.Range("E23").Value = "2.06%"
$ws.Range("D24").Value = "0.001248"
$ws.Range("E24").Value = "3.77%"
$ws.Range("D25").Value = "0.004354"
$ws.Range("E25").Value = "1.94%"
$ws.Range("D26").Value = "0.0001187"
$ws.Range("E26").Value = "-15.45%"
$ws.Range("D27").Value = "0.0003980"
$ws.Range("E27").Value = "36.92%"
$ws.Range("D39").Value = "0.02521"
$ws.Range("E39").Value = "18.54%"
$ws.Range("D40").Value = "0.05247"
$ws.Range("E40").Value = "6.29%"
$ws.Range("D41").Value = "0.007298"
$ws.Range("E41").Value = "-3.85%"
$ws.Range("D42").Value = "0.1436"
$ws.Range("E42").Value = "6.44%"
$ws.Range("D43").Value = "0.008839"
$ws.Range("E43").Value = "3.85%"
$ws.Range("D44").Value = "0.002181"
$ws.Range("E44").Value = "9.48%"
$ws.Range("D45").Value = "0.009634"
$ws.Range("E45").Value = "11.87%"
$ws.Range("D46").Value = "0.00006718"
$ws.Range("E46").Value = "1.89%"
$ws.Range("D47").Value = "0.00000000748"
$ws.Range("E47").Value = "-0.48%"
$ws.Range("D48").Value = "0.002867"
$ws.Range("E48").Value = "-13.14%"
$ws.Range("D49").Value = "0.001796"
$ws.Range("E49").Value = "24.15%"
$ws.Range("D50").Value = "0.00002095"
$ws.Range("E50").Value = "-0.48%"
$ws.Range("D51").Value = "0.0001995"
$ws.Range("E51").Value = "-0.48%"

# Restore original (default) style so no stray number-format/style
# attribute is left behind on these cells.
$rng.Style = "Normal"
